$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New s_val data (regenerated to filter save games) for rows 2-22.
# Columns: B=TB, C=d2S, D=K, E=IP, G=sum (F=Win is unchanged)
$data = @{
    2 = @{ B = 1.445647641019636; C = 1.626987699542094; D = 3.223369029078222; E = 13.86384647080068; G = 20.15985084044064 }
    3 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    4 = @{ B = 1.445647641019636; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    5 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    6 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    7 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 3.223369029078222; E = 0.5333859586016987; G = 8.656069925401464 }
    8 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 3.223369029078222; E = 0.5333859586016987; G = 8.656069925401464 }
    9 = @{ B = 0.2881169905109251; C = 0.3048912486333797; D = 3.223369029078222; E = 0.5333859586016987; G = 4.349763226824225 }
    10 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 3.223369029078222; E = 0.5333859586016987; G = 8.656069925401464 }
    11 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
    12 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    13 = @{ B = 0.6545652718822623; C = 0.3048912486333797; D = 3.223369029078222; E = 0.5333859586016987; G = 4.716211508195562 }
    14 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.1496068669990043; E = 13.86384647080068; G = 18.91276827552123 }
    15 = @{ B = 1.445647641019636; C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 3.755628166162433 }
    16 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    17 = @{ B = 0.1169995834814548; C = 0.3048912486333797; D = 0.7210945179870265; E = 0.5333859586016987; G = 1.67637130870356 }
    18 = @{ B = 1.445647641019636; C = 0.3048912486333797; D = 3.223369029078222; E = 0.5333859586016987; G = 5.507293877332936 }
    19 = @{ B = 1.445647641019636; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    20 = @{ B = 0.6545652718822623; C = 1.626987699542094; D = 18.71679738969934; E = 0.5333859586016987; G = 21.53173631972539 }
    21 = @{ B = 1.445647641019636; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    22 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
}

foreach ($rowNum in $data.Keys) {
    $row = $data[$rowNum]
    $ws.Range("B$rowNum").Value = $row.B
    $ws.Range("C$rowNum").Value = $row.C
    $ws.Range("D$rowNum").Value = $row.D
    $ws.Range("E$rowNum").Value = $row.E
    $ws.Range("G$rowNum").Value = $row.G
}
